$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.654.13"
$ws.Range("E2").Value = "  +4.75%  "

$ws.Range("D3").Value = "2.491.10"
$ws.Range("E3").Value = "  +2.53%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'323.47"
$ws.Range("E5").Value = "  +1.60%  "

$ws.Range("E6").Value = "  +2.22%  "

$ws.Range("E7").Value = "  +1.80%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "'0.544"
$ws.Range("E9").Value = "  +2.34%  "

$ws.Range("D10").Value = "'38.37"
$ws.Range("E10").Value = "  +7.45%  "

$ws.Range("E11").Value = "  +1.46%  "

$ws.Range("E12").Value = "  +1.14%  "

$ws.Range("D13").Value = "'18.32"
$ws.Range("E13").Value = "  +0.96%  "

$ws.Range("E14").Value = "  +1.53%  "

$ws.Range("D15").Value = "2.877.71"
$ws.Range("E15").Value = "  +2.47%  "

$ws.Range("D16").Value = "2.503.53"
$ws.Range("E16").Value = "  +3.27%  "

$ws.Range("D17").Value = "'0.846"
$ws.Range("E17").Value = "  +0.21%  "

$ws.Range("D18").Value = "47.505.83"
$ws.Range("E18").Value = "  +4.75%  "

$ws.Range("D19").Value = "'12.75"
$ws.Range("E19").Value = "  +3.93%  "

$ws.Range("D20").Value = "'6.59"
$ws.Range("E20").Value = "  +3.68%  "

$ws.Range("D21").Value = "0.0₃0938"
$ws.Range("E21").Value = "  +1.69%  "

$ws.Range("D22").Value = "'70.76"
$ws.Range("E22").Value = "  +2.72%  "

$ws.Range("D23").Value = "'251.31"
$ws.Range("E23").Value = "  +2.66%  "

$ws.Range("E24").Value = "  +5.75%  "

$ws.Range("D25").Value = "'2.58"
$ws.Range("E25").Value = "  +3.14%  "

$ws.Range("D26").Value = "'26.17"
$ws.Range("E26").Value = "  +2.30%  "

$ws.Range("E27").Value = "  +0.04%  "

$ws.Range("E28").Value = "  +4.66%  "

$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "'35.17"
$ws.Range("E29").Value = "  +6.79%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.08"
$ws.Range("E30").Value = "  -4.54%  "

$ws.Range("E31").Value = "  +7.47%  "

$ws.Range("E32").Value = "  +0.47%  "

$ws.Range("E33").Value = "  -2.39%  "

$ws.Range("D34").Value = "'5.36"
$ws.Range("E34").Value = "  +2.62%  "

$ws.Range("D35").Value = "'0.0784"
$ws.Range("E35").Value = "  +1.93%  "

$ws.Range("E36").Value = "  +0.19%  "

$ws.Range("D37").Value = "'1.97"
$ws.Range("E37").Value = "  +4.78%  "

$ws.Range("D38").Value = "'4.65"
$ws.Range("E38").Value = "  +4.25%  "

$ws.Range("E39").Value = "  +4.50%  "

$ws.Range("E40").Value = "  +1.85%  "

$ws.Range("E41").Value = "  +1.90%  "

$ws.Range("D42").Value = "'122.24"
$ws.Range("E42").Value = "  -2.46%  "

$ws.Range("D43").Value = "'21.17"
$ws.Range("E43").Value = "  +3.87%  "

$ws.Range("E44").Value = "  +2.58%  "

$ws.Range("D45").Value = "1.966.95"

$ws.Range("E46").Value = "  +1.59%  "

$ws.Range("E47").Value = "  -0.23%  "

$ws.Range("D48").Value = "'9.21"
$ws.Range("E48").Value = "  +1.05%  "

$ws.Range("E49").Value = "  -1.05%  "

$ws.Range("D50").Value = "'5.29"
$ws.Range("E50").Value = "  +11.66%  "

$ws.Range("D51").Value = "'79.50"
$ws.Range("E51").Value = "  +3.46%  "
